# Reorder the "Recorded By" (column G) values on the "Session Analysis Results"
# sheet so that each comma-separated list of recorders is written back in a
# different (swapped) order, matching the upstream sync of the source report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of old value -> new value, as observed in the diff.
$map = @{
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
